$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.670.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.346.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.83%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.46%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.54%  "

$ws.Range("E7").Value = "  -2.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.616"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.31"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.51%  "

$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -8.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.700.43"
$ws.Range("D16").Style = "Normal"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.377.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.648.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.68"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.91%  "

$ws.Range("E20").Value = "  -2.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "77.10"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.39%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.61"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.65"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.19%  "

$ws.Range("E24").Value = "  -5.27%  "

$ws.Range("E25").Value = "  -5.27%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -4.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.20%  "

$ws.Range("E29").Value = "  +3.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "174.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.48%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.62%  "

$ws.Range("E32").Value = "  -4.62%  "

$ws.Range("E33").Value = "  +3.00%  "

$ws.Range("E34").Value = "  -8.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.123"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +12.45%  "

$ws.Range("E36").Value = "  -3.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.59"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.89%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.76"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.99%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "71.73"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.235"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("E43").Value = "  -8.18%  "

$ws.Range("E44").Value = "  -0.07%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "114.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "11.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "5.47"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "

$ws.Range("E49").Value = "  -5.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.39%  "

$ws.Range("E51").Value = "  -1.79%  "
